# SwaadSutra Daily Orders update — 2026-01-19T05:41:03.627Z
#
# A new order (#16, same customer/items as the previous order but placed
# two minutes later, with a later collection date and no special notes)
# comes in. It is inserted as the new row 2 (most-recent-first listing);
# the former row 2 (order #15) is pushed down to row 3, unchanged. The
# Summary and Items Breakdown sheets are updated to reflect the new
# (now-doubled) totals.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet: Daily Orders
# ---------------------------------------------------------------------
$orders = $wb.Worksheets.Item("Daily Orders")

# Duplicate row 2 (order #15) down into row 3 via Copy so every cell's
# original type (text vs number) is preserved exactly - a plain .Value
# assignment would let Excel "helpfully" reinterpret things like the
# phone number or the collection date as numbers/dates.
$orders.Range("A2:N2").Copy($orders.Range("A3:N3"))

# Turn row 2 into the new order (#16): new order id, later timestamp,
# later collection date, and no notes this time.
$orders.Range("A2").Value = 16
$orders.Range("B2").Value = "2026-01-19 05:41"
$orders.Range("J2").NumberFormat = "@"
$orders.Range("J2").Value = "2026-01-21"
$orders.Range("L2").Value = ""

# ---------------------------------------------------------------------
# Sheet: Summary
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("A2").Value = 2     # Total Orders
$summary.Range("B2").Value = 2     # New
$summary.Range("G2").Value = 210   # Total Revenue

# ---------------------------------------------------------------------
# Sheet: Items Breakdown
# ---------------------------------------------------------------------
$items = $wb.Worksheets.Item("Items Breakdown")
$items.Range("B2").Value = 10   # Wheat Chapati - Quantity Ordered
$items.Range("C2").Value = 150  # Wheat Chapati - Revenue
$items.Range("B3").Value = 2    # 1 Plate Bhaji - Quantity Ordered
$items.Range("C3").Value = 60   # 1 Plate Bhaji - Revenue
